$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -10
$ws.Range("F3").Value = -4
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = 2
